{"js": "// 20.04 Arduino + Android MQ7 updates\n// Append, at the very end of the document (after the existing ARDUINO /\n// HM-10 Bluetooth hyperlink paragraph), a blank paragraph followed by a\n// new paragraph with the MQ-7 sensor code link:\n//   <blank paragraph>\n//   <TAB>Cod pentru MQ-7:  <hyperlink to the MQ-7 tutorial>\n\nconst url = \"https://www.teachmemicro.com/use-mq-7-carbon-monoxide-sensor/\";\nconst label = \"\\tCod pentru MQ-7 : \";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Last paragraph currently in the document (the one ending with the\n// HM-10 Bluetooth hyperlink).\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// 1) Insert a blank paragraph right after it.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// 2) Insert a second new paragraph holding the label text.\nconst codeParagraph = blankParagraph.insertParagraph(label, \"After\");\nawait context.sync();\n\n// 3) Append the URL text at the end of that paragraph, then turn just\n//    that range into a hyperlink.\nconst endRange = codeParagraph.getRange(\"End\");\nconst linkRange = endRange.insertText(url, \"End\");\nlinkRange.hyperlink = url;\nawait context.sync();\n", "ps1": "# 20.04 Arduino + Android MQ7 updates\n# Append, at the very end of the document (after the existing ARDUINO /\n# HM-10 Bluetooth hyperlink paragraph), a blank paragraph followed by a\n# new paragraph with the MQ-7 sensor code link:\n#   <blank paragraph>\n#   <TAB>Cod pentru MQ-7:  <hyperlink to the MQ-7 tutorial>\n\n$d = $word.ActiveDocument\n\n$url = \"https://www.teachmemicro.com/use-mq-7-carbon-monoxide-sensor/\"\n$label = [char]9 + \"Cod pentru MQ-7 : \"\n\n# 1) Blank paragraph right after the last paragraph in the document.\n$tail = $d.Range($d.Content.End - 1, $d.Content.End - 1)\n$tail.InsertParagraphAfter()\n\n# 2) A second new paragraph that will hold the label text + hyperlink.\n$tail = $d.Range($d.Content.End - 1, $d.Content.End - 1)\n$tail.InsertParagraphAfter()\n\n# 3) Insert the label text followed immediately by the raw URL text into\n#    that new (currently empty) paragraph.\n$insStart = $d.Content.End - 1\n$ins = $d.Range($insStart, $insStart)\n$ins.InsertAfter($label + $url)\n\n# 4) Locate the URL text we just typed and convert it into a real\n#    hyperlink (this keeps the label text as plain runs before the\n#    hyperlink, matching how Word builds the paragraph).\n$urlRange = $d.Range($insStart, $d.Content.End - 1)\n$urlRange.Find.ClearFormatting()\n[void]$urlRange.Find.Execute($url)\n$d.Hyperlinks.Add($urlRange, $url)\n"}
